# Update the "Generate Report for Handback" timestamps across the sheets.
$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" (G2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-03 17:12:52"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-03 17:12:47"
$wsZhCn.Range("K2").Value = "2016-09-03 17:13:09"

# de-de sheet: "Correspond Handback DateTime" (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-03 17:13:17"
